# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# 87125e70... and f9fde33b... items have been handed back (in addition to
# already being handed off), for both the zh-cn and de-de target languages:
#   - "Status" on the Overview sheet flips from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - The "Latest Target File" (I) and "Latest Handback File" (J) columns on
#     the zh-cn / de-de sheets get populated (were blank), and a hyperlink is
#     added on the target-file cell, matching the existing hyperlink already
#     present on column A.
#   - The "Latest Handback DateTime" (K) column is stamped with the handback
#     timestamp.
#   - A handful of columns are widened to comfortably fit the newly
#     populated long file names.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/70c715776eb314f56c976c5624cd64c598547645/e2e/"

$file1Id   = "87125e70-4812-4d7a-bab9-591f8a17caf5"
$file2Id   = "f9fde33b-ca70-47b4-998d-c05cc45437ce"
$file1Md   = "$file1Id.md"
$file2Md   = "$file2Id.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status goes from "Ready for handoff" to "Handed back"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# Widen the now-longer status columns.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): fill in handback info
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-08-13 17:24:17" },
    @{ Name = "de-de"; HandbackTime = "2016-08-13 17:24:27" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # "Latest Handback File" (J) is the very same xliff that was just handed
    # off, i.e. it mirrors "Latest Handoff File" (G) for that row.
    $handback1 = $ws.Range("G2").Text
    $handback2 = $ws.Range("G3").Text

    # Row 2 -> file 1 (87125e70...)
    $ws.Range("I2").Value = $file1Md
    $ws.Range("J2").Value = $handback1
    $ws.Range("K2").Value = $lang.HandbackTime

    # Row 3 -> file 2 (f9fde33b...)
    $ws.Range("I3").Value = $file2Md
    $ws.Range("J3").Value = $handback2
    $ws.Range("K3").Value = $lang.HandbackTime

    # Give the newly-populated "Latest Target File" cells the same visual
    # style (the built-in HyperLink style) as the existing link in column A.
    $ws.Range("I2").Style = "Hyperlink"
    $ws.Range("I3").Style = "Hyperlink"

    # Add the actual hyperlinks, mirroring the ones on column A.
    $ws.Hyperlinks.Add($ws.Range("I2"), ($repoBase + $file1Md), [Type]::Missing, [Type]::Missing, $file1Md)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($repoBase + $file2Md), [Type]::Missing, [Type]::Missing, $file2Md)

    # Widen columns so the long file names/status text are fully visible.
    $ws.Columns.Item(3).ColumnWidth = 29.14   # Status
    $ws.Columns.Item(9).ColumnWidth = 39.17   # Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.17  # Latest Handback File
}
